$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new blank row at row 2. This shifts the existing rows
#    2..7 down to rows 3..8 (their content is untouched).
# ------------------------------------------------------------------
$ws.Rows("2:2").Insert()

# The freshly inserted row 2 inherits the bold/centered header
# formatting from row 1 above it - strip that back to the plain
# (unstyled) look the data rows use.
$ws.Range("A2:Q2").ClearFormats()

# ------------------------------------------------------------------
# 2) The row that used to be row 7 (Milena Alves Barboza ... /
#    PREVALENCIA E OBITOS ...) is now sitting at row 8 after the
#    insert above. Copy its values up into row 2 verbatim.
# ------------------------------------------------------------------
$ws.Range("A2:Q2").NumberFormat = "@"
$ws.Range("A2:Q2").Value = $ws.Range("A8:Q8").Value()
$ws.Range("A2:Q2").ClearFormats()

# ------------------------------------------------------------------
# 3) Row 8 now holds a duplicate of that same data - overwrite it
#    with the brand-new "Fernanda Coutinho ..." / "Gait Patterns ..."
#    record that belongs there.
# ------------------------------------------------------------------
$ws.Range("A8:Q8").NumberFormat = "@"

$ws.Cells.Item(8,1).Value = "Fernanda Coutinho, José Ribeiro Ferreira, Nuno Nogueira"
$ws.Cells.Item(8,2).Value = "; ; "
$ws.Cells.Item(8,3).Value = "https://openalex.org/W4391913442"
$ws.Cells.Item(8,4).Value = "Gait Patterns Analysis Using Pressure Plataforms"
$ws.Cells.Item(8,5).Value = "2023-01-01"
$ws.Cells.Item(8,6).Value = "Lecture notes in bioengineering"
$ws.Cells.Item(8,7).Value = "Springer International Publishing"
$ws.Cells.Item(8,8).Value = "https://doi.org/10.1007/978-3-031-47790-4_56"
$ws.Cells.Item(8,9).Value = "N/A"
$ws.Cells.Item(8,10).Value = "N/A"
$ws.Cells.Item(8,11).Value = "closed"
$ws.Cells.Item(8,12).Value = "en"
$ws.Cells.Item(8,13).Value = "0"
$ws.Cells.Item(8,14).Value = "2023"
$ws.Cells.Item(8,15).Value = "NA"
$ws.Cells.Item(8,16).Value = "https://doi.org/10.1007/978-3-031-47790-4_56"
$ws.Cells.Item(8,17).Value = "book-chapter"

$ws.Range("A8:Q8").ClearFormats()

# ------------------------------------------------------------------
# 4) The record that used to be row 4 (the "Rotary-Cage Valve" paper)
#    now lives at row 5. Its host_organization (column G) is
#    corrected from "N/A" to the actual publisher name.
# ------------------------------------------------------------------
$ws.Range("G5").Value = "Springer Science+Business Media"
